{"js": "const replacements = [\n  [\"658\u00f79=73, 1\", \"669\u00f76=111, 3\"],\n  [\"532\u00f73=177, 1\", \"393\u00f77=56, 1\"],\n  [\"941\u00f72=470, 1\", \"267\u00f73=89, 0\"],\n  [\"577\u00f77=82, 3\", \"703\u00f76=117, 1\"],\n  [\"938\u00f78=117, 2\", \"467\u00f75=93, 2\"],\n  [\"372\u00f75=74, 2\", \"988\u00f77=141, 1\"],\n  [\"814\u00f77=116, 2\", \"578\u00f79=64, 2\"],\n  [\"802\u00f78=100, 2\", \"945\u00f75=189, 0\"],\n  [\"440\u00f72=220, 0\", \"881\u00f75=176, 1\"],\n  [\"148\u00f73=49, 1\", \"651\u00f73=217, 0\"],\n  [\"246\u00f78=30, 6\", \"287\u00f79=31, 8\"],\n  [\"340\u00f79=37, 7\", \"758\u00f72=379, 0\"],\n  [\"557\u00f73=185, 2\", \"186\u00f78=23, 2\"],\n  [\"811\u00f74=202, 3\", \"743\u00f72=371, 1\"],\n  [\"639\u00f73=213, 0\", \"236\u00f72=118, 0\"],\n  [\"663\u00f79=73, 6\", \"732\u00f74=183, 0\"],\n  [\"494\u00f78=61, 6\", \"686\u00f79=76, 2\"],\n  [\"854\u00f72=427, 0\", \"960\u00f78=120, 0\"],\n  [\"834\u00f75=166, 4\", \"231\u00f78=28, 7\"],\n  [\"584\u00f77=83, 3\", \"849\u00f79=94, 3\"],\n  [\"152\u00f72=76, 0\", \"483\u00f73=161, 0\"],\n  [\"964\u00f73=321, 1\", \"195\u00f79=21, 6\"],\n  [\"951\u00f74=237, 3\", \"769\u00f78=96, 1\"],\n  [\"277\u00f78=34, 5\", \"584\u00f72=292, 0\"],\n  [\"362\u00f77=51, 5\", \"189\u00f77=27, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\nreturn \"ok\";\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"658\u00f79=73, 1\", \"669\u00f76=111, 3\"),\n    @(\"532\u00f73=177, 1\", \"393\u00f77=56, 1\"),\n    @(\"941\u00f72=470, 1\", \"267\u00f73=89, 0\"),\n    @(\"577\u00f77=82, 3\", \"703\u00f76=117, 1\"),\n    @(\"938\u00f78=117, 2\", \"467\u00f75=93, 2\"),\n    @(\"372\u00f75=74, 2\", \"988\u00f77=141, 1\"),\n    @(\"814\u00f77=116, 2\", \"578\u00f79=64, 2\"),\n    @(\"802\u00f78=100, 2\", \"945\u00f75=189, 0\"),\n    @(\"440\u00f72=220, 0\", \"881\u00f75=176, 1\"),\n    @(\"148\u00f73=49, 1\", \"651\u00f73=217, 0\"),\n    @(\"246\u00f78=30, 6\", \"287\u00f79=31, 8\"),\n    @(\"340\u00f79=37, 7\", \"758\u00f72=379, 0\"),\n    @(\"557\u00f73=185, 2\", \"186\u00f78=23, 2\"),\n    @(\"811\u00f74=202, 3\", \"743\u00f72=371, 1\"),\n    @(\"639\u00f73=213, 0\", \"236\u00f72=118, 0\"),\n    @(\"663\u00f79=73, 6\", \"732\u00f74=183, 0\"),\n    @(\"494\u00f78=61, 6\", \"686\u00f79=76, 2\"),\n    @(\"854\u00f72=427, 0\", \"960\u00f78=120, 0\"),\n    @(\"834\u00f75=166, 4\", \"231\u00f78=28, 7\"),\n    @(\"584\u00f77=83, 3\", \"849\u00f79=94, 3\"),\n    @(\"152\u00f72=76, 0\", \"483\u00f73=161, 0\"),\n    @(\"964\u00f73=321, 1\", \"195\u00f79=21, 6\"),\n    @(\"951\u00f74=237, 3\", \"769\u00f78=96, 1\"),\n    @(\"277\u00f78=34, 5\", \"584\u00f72=292, 0\"),\n    @(\"362\u00f77=51, 5\", \"189\u00f77=27, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $ok = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Replacement not found: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
